# Auto-generated edit script applying scheduled market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across several
# crafting-class sheets, matching the upstream commit's data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 569.8
$ws.Range("I5").Value = 683
$ws.Range("K5").Value = 683
$ws.Range("M5").Value = -568

$ws.Range("H15").Value = 1748.8823
$ws.Range("I15").Value = 1748.8823
$ws.Range("K15").Value = 5246.6469
$ws.Range("M15").Value = -5077.6469

$ws.Range("H64").Value = 3942.2424
$ws.Range("I64").Value = 3438.111
$ws.Range("J64").Value = 4131.2915
$ws.Range("K64").Value = 3438.111
$ws.Range("L64").Value = 4131.2915
$ws.Range("M64").Value = -3190.111
$ws.Range("N64").Value = -4627.2915

$ws.Range("H67").Value = 3942.2424
$ws.Range("I67").Value = 3438.111
$ws.Range("J67").Value = 4131.2915
$ws.Range("K67").Value = 3438.111
$ws.Range("L67").Value = 4131.2915
$ws.Range("M67").Value = -2580.111
$ws.Range("N67").Value = -5847.2915

$ws.Range("H137").Value = 3574110.8
$ws.Range("I137").Value = 5001926
$ws.Range("J137").Value = 4572.375
$ws.Range("K137").Value = 15005778
$ws.Range("L137").Value = 13717.125
$ws.Range("M137").Value = -15003228
$ws.Range("N137").Value = -18817.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11945962
$ws.Range("I32").Value = 15201823
$ws.Range("J32").Value = 7802.1665
$ws.Range("K32").Value = 15201823
$ws.Range("L32").Value = 7802.1665
$ws.Range("M32").Value = -15201536
$ws.Range("N32").Value = -8376.166499999999

$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992

$ws.Range("H122").Value = 1091.8387
$ws.Range("I122").Value = 1087.3448
$ws.Range("J122").Value = 1157
$ws.Range("K122").Value = 3262.0344
$ws.Range("L122").Value = 3471
$ws.Range("M122").Value = -812.0344000000005
$ws.Range("N122").Value = -8371

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 60000
$ws.Range("J55").Value = 60000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3135.8223
$ws.Range("I31").Value = 2044.8
$ws.Range("J31").Value = 3200
$ws.Range("K31").Value = 2044.8
$ws.Range("L31").Value = 3200
$ws.Range("M31").Value = -1749.8
$ws.Range("N31").Value = -3790

$ws.Range("H34").Value = 3135.8223
$ws.Range("I34").Value = 2044.8
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 2044.8
$ws.Range("L34").Value = 3200
$ws.Range("M34").Value = -1842.8
$ws.Range("N34").Value = -3604

$ws.Range("H132").Value = 50949.43
$ws.Range("I132").Value = 2730
$ws.Range("J132").Value = 171498
$ws.Range("K132").Value = 8190
$ws.Range("L132").Value = 514494
$ws.Range("M132").Value = -5660
$ws.Range("N132").Value = -519554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7223433
$ws.Range("I4").Value = 3000779
$ws.Range("J4").Value = 12501750
$ws.Range("K4").Value = 9002337
$ws.Range("L4").Value = 37505250
$ws.Range("M4").Value = -9002225
$ws.Range("N4").Value = -37505474

$ws.Range("H112").Value = 17547118
$ws.Range("I112").Value = 1933.3334
$ws.Range("J112").Value = 20836840
$ws.Range("K112").Value = 5800.0002
$ws.Range("L112").Value = 62510520
$ws.Range("M112").Value = -4692.0002
$ws.Range("N112").Value = -62512736

$ws.Range("H132").Value = 4000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

$ws.Range("H133").Value = 3347.1428
$ws.Range("I133").Value = 3347.1428
$ws.Range("K133").Value = 10041.4284
$ws.Range("M133").Value = -4981.428400000001

$ws.Range("H134").Value = 2419.5881
$ws.Range("I134").Value = 1511.1666
$ws.Range("J134").Value = 4599.8
$ws.Range("K134").Value = 4533.4998
$ws.Range("L134").Value = 13799.4
$ws.Range("M134").Value = 536.5002000000004
$ws.Range("N134").Value = -23939.4

$ws.Range("H137").Value = 32618.316
$ws.Range("J137").Value = 121818
$ws.Range("L137").Value = 365454
$ws.Range("N137").Value = -375654

$ws.Range("H139").Value = 337584.3
$ws.Range("I139").Value = 627156.0600000001
$ws.Range("J139").Value = 6645.143
$ws.Range("K139").Value = 1881468.18
$ws.Range("L139").Value = 19935.429
$ws.Range("M139").Value = -1876328.18
$ws.Range("N139").Value = -30215.429

$ws.Range("H141").Value = 3286.3635
$ws.Range("I141").Value = 3286.3635
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9859.0905
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4679.0905
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4624
$ws.Range("I122").Value = 4498.6665
$ws.Range("K122").Value = 13495.9995
$ws.Range("M122").Value = -11045.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1430.8096
$ws.Range("I68").Value = 1363.7222
$ws.Range("J68").Value = 1833.3334
$ws.Range("K68").Value = 1363.7222
$ws.Range("L68").Value = 1833.3334
$ws.Range("M68").Value = -614.7221999999999
$ws.Range("N68").Value = -3331.3334

$ws.Range("H71").Value = 1430.8096
$ws.Range("I71").Value = 1363.7222
$ws.Range("J71").Value = 1833.3334
$ws.Range("K71").Value = 6818.611
$ws.Range("L71").Value = 9166.666999999999
$ws.Range("M71").Value = -3074.611
$ws.Range("N71").Value = -16654.667

$ws.Range("H122").Value = 2613.45
$ws.Range("I122").Value = 2529.647
$ws.Range("J122").Value = 3088.3333
$ws.Range("K122").Value = 7588.941
$ws.Range("L122").Value = 9264.999899999999
$ws.Range("M122").Value = -5138.941
$ws.Range("N122").Value = -14164.9999

$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

$ws.Range("H136").Value = 722125
$ws.Range("I136").Value = 504250
$ws.Range("J136").Value = 940000
$ws.Range("K136").Value = 1512750
$ws.Range("L136").Value = 2820000
$ws.Range("M136").Value = -1510200
$ws.Range("N136").Value = -2825100

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

